# Fruta / hortaliza, semanal
#
# The weekly refresh re-sorts the price records (rows 2-25) of the
# "Camote" sheet. Column A/B/C/E/F/G/H/I/R are constant for every record
# in this sub-sheet, so the only visible effect of the re-sort is that
# the data carried in columns D (Fecha), J (Volumen), K/L/M (precios),
# N (Unidad de comercializacion), O (Origen), P (Precio $/Kg) and
# Q (Kg o Unidades) moves from one row to another.
#
# Capture each source row's values first (so that we don't clobber a
# value before it has been copied elsewhere), then write every
# destination row from its captured source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (where the data used to live)
$rowMap = @{
    2  = 18
    3  = 19
    4  = 20
    5  = 9
    6  = 2
    7  = 10
    8  = 21
    9  = 4
    10 = 3
    11 = 7
    12 = 8
    13 = 12
    14 = 11
    15 = 15
    16 = 22
    17 = 14
    18 = 25
    19 = 5
    20 = 16
    21 = 17
    22 = 13
    23 = 23
    24 = 24
    25 = 6
}

# Columns that travel together with a record when rows get re-sorted.
$cols = @(4, 10, 11, 12, 13, 14, 15, 16, 17)

# Snapshot every source row's values before writing anything.
# (.Value2 is used instead of .Value for reading raw numbers/strings.)
$snapshot = @{}
foreach ($row in 2..25) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowValues
}

# Now write each destination row using the snapshot of its source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value2 = $srcValues[$col]
    }
}
